$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update the repayment-strategy value on the ProductLoanInput sheet: the old
# "RBI (India)" scenario is replaced with the new periodic/upfront ordering
# scenario, with a left/top-aligned style for the longer text.
$ws1.Range("B17").Value = "Penalties, Fees, Interest, Principal order"
$ws1.Range("B17").HorizontalAlignment = -4131
$ws1.Range("B17").VerticalAlignment = -4160

# Make ProductLoanInput the active/selected tab (it was ProductLoanOutput
# before), with the selection resting on the cell we just edited.
[void]$ws1.Activate()
[void]$ws1.Range("B17").Select()
